$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 115
$ws.Cells.Item(115, 1).Value = 109125567
$ws.Cells.Item(115, 2).Value = "2021.09.16 11:48:46"
$ws.Cells.Item(115, 3).Value = "buy"
$ws.Cells.Item(115, 4).Value = 0.02
$ws.Cells.Item(115, 5).Value = "us.100.."
$ws.Cells.Item(115, 6).Value = 15453.67
$ws.Cells.Item(115, 7).Value = 15453.35
$ws.Cells.Item(115, 8).Value = 15511
$ws.Cells.Item(115, 9).Value = "2021.09.16 13:35:07"
$ws.Cells.Item(115, 10).Value = 15453.35
$ws.Cells.Item(115, 11).Value = -0.64
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 13).Value = 0
$ws.Cells.Item(115, 14).Value = -0.5

# Row 116
$ws.Cells.Item(116, 1).Value = 109122880
$ws.Cells.Item(116, 2).Value = "2021.09.16 06:50:02"
$ws.Cells.Item(116, 3).Value = "buy"
$ws.Cells.Item(116, 4).Value = 0.02
$ws.Cells.Item(116, 5).Value = "us.100.."
$ws.Cells.Item(116, 6).Value = 15487.69
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 15512
$ws.Cells.Item(116, 9).Value = "2021.09.16 21:07:24"
$ws.Cells.Item(116, 10).Value = 15512
$ws.Cells.Item(116, 11).Value = -0.64
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = 0
$ws.Cells.Item(116, 14).Value = 37.76

# Row 117
$ws.Cells.Item(117, 1).Value = 109122889
$ws.Cells.Item(117, 2).Value = "2021.09.16 06:56:13"
$ws.Cells.Item(117, 3).Value = "buy"
$ws.Cells.Item(117, 4).Value = 0.02
$ws.Cells.Item(117, 5).Value = "us.100.."
$ws.Cells.Item(117, 6).Value = 15485.67
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 15512
$ws.Cells.Item(117, 9).Value = "2021.09.16 21:07:24"
$ws.Cells.Item(117, 10).Value = 15512
$ws.Cells.Item(117, 11).Value = -0.64
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 13).Value = 0
$ws.Cells.Item(117, 14).Value = 40.9

# Row 118
$ws.Cells.Item(118, 1).Value = 109122983
$ws.Cells.Item(118, 2).Value = "2021.09.16 07:24:09"
$ws.Cells.Item(118, 3).Value = "buy"
$ws.Cells.Item(118, 4).Value = 0.02
$ws.Cells.Item(118, 5).Value = "us.100.."
$ws.Cells.Item(118, 6).Value = 15481.06
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 15512
$ws.Cells.Item(118, 9).Value = "2021.09.16 21:07:24"
$ws.Cells.Item(118, 10).Value = 15512
$ws.Cells.Item(118, 11).Value = -0.64
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 13).Value = 0
$ws.Cells.Item(118, 14).Value = 48.06

# Row 119
$ws.Cells.Item(119, 1).Value = 109123618
$ws.Cells.Item(119, 2).Value = "2021.09.16 08:52:51"
$ws.Cells.Item(119, 3).Value = "buy"
$ws.Cells.Item(119, 4).Value = 0.02
$ws.Cells.Item(119, 5).Value = "us.100.."
$ws.Cells.Item(119, 6).Value = 15473.06
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 15512
$ws.Cells.Item(119, 9).Value = "2021.09.16 21:07:24"
$ws.Cells.Item(119, 10).Value = 15512
$ws.Cells.Item(119, 11).Value = -0.64
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 13).Value = 0
$ws.Cells.Item(119, 14).Value = 60.48

# Row 120
$ws.Cells.Item(120, 1).Value = 109124227
$ws.Cells.Item(120, 2).Value = "2021.09.16 09:20:18"
$ws.Cells.Item(120, 3).Value = "buy"
$ws.Cells.Item(120, 4).Value = 0.02
$ws.Cells.Item(120, 5).Value = "us.100.."
$ws.Cells.Item(120, 6).Value = 15495.31
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 15512
$ws.Cells.Item(120, 9).Value = "2021.09.16 21:07:24"
$ws.Cells.Item(120, 10).Value = 15512
$ws.Cells.Item(120, 11).Value = -0.64
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(120, 13).Value = 0
$ws.Cells.Item(120, 14).Value = 25.92

# Row 121
$ws.Cells.Item(121, 1).Value = 109126351
$ws.Cells.Item(121, 2).Value = "2021.09.16 13:45:09"
$ws.Cells.Item(121, 3).Value = "buy"
$ws.Cells.Item(121, 4).Value = 0.02
$ws.Cells.Item(121, 5).Value = "us.100.."
$ws.Cells.Item(121, 6).Value = 15454.06
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 15512
$ws.Cells.Item(121, 9).Value = "2021.09.16 21:07:24"
$ws.Cells.Item(121, 10).Value = 15512
$ws.Cells.Item(121, 11).Value = -0.64
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 13).Value = 0
$ws.Cells.Item(121, 14).Value = 89.99

# Update selection to match new active cell
$ws.Range("P114").Select() | Out-Null
